{"js": "// Replace the trailing empty paragraph with one containing \"Be real man!\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst last = items[items.length - 1];\n\nif (last.text === \"\") {\n  // Fill the existing empty trailing paragraph with the new text.\n  last.insertText(\"Be real man!\", Word.InsertLocation.replace);\n} else {\n  // Fallback: if there's no empty trailing paragraph, add a new one.\n  body.insertParagraph(\"Be real man!\", Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# Replace the trailing empty paragraph with one containing \"Be real man!\"\n$d = $word.ActiveDocument\n$p = $d.Paragraphs.Last\n\n# Paragraph.Range.Text includes the trailing paragraph-mark (chr 13), so an\n# \"empty\" paragraph's text is just that single character.\nif ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"\") {\n    $p.Range.InsertBefore(\"Be real man!\")\n} else {\n    $d.Content.InsertParagraphAfter()\n    $d.Paragraphs.Last.Range.Text = \"Be real man!\"\n}\n"}
